$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift every timestamp in column A (rows 2-97) forward by 4 days,
# preserving the time-of-day fraction exactly.
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 4
}

# Updated forecast values in column B for the rows that changed.
$bUpdates = @(
    @{Row=18; Value=6},
    @{Row=19; Value=6},
    @{Row=20; Value=8},
    @{Row=21; Value=12},
    @{Row=22; Value=113},
    @{Row=23; Value=125},
    @{Row=24; Value=142},
    @{Row=25; Value=164},
    @{Row=26; Value=464},
    @{Row=27; Value=493},
    @{Row=28; Value=527},
    @{Row=29; Value=567},
    @{Row=30; Value=1071},
    @{Row=31; Value=1113},
    @{Row=32; Value=1160},
    @{Row=33; Value=1208},
    @{Row=34; Value=1589},
    @{Row=35; Value=1631},
    @{Row=36; Value=1673},
    @{Row=37; Value=1714},
    @{Row=38; Value=1959},
    @{Row=39; Value=2004},
    @{Row=40; Value=2036},
    @{Row=41; Value=2062},
    @{Row=42; Value=2199},
    @{Row=43; Value=2214},
    @{Row=44; Value=2227},
    @{Row=45; Value=2239},
    @{Row=46; Value=2293},
    @{Row=47; Value=2299},
    @{Row=48; Value=2304},
    @{Row=49; Value=2303},
    @{Row=50; Value=2288},
    @{Row=51; Value=2285},
    @{Row=52; Value=2280},
    @{Row=53; Value=2272},
    @{Row=54; Value=2160},
    @{Row=55; Value=2149},
    @{Row=56; Value=2136},
    @{Row=57; Value=2124},
    @{Row=58; Value=2004},
    @{Row=59; Value=1981},
    @{Row=60; Value=1958},
    @{Row=61; Value=1932},
    @{Row=62; Value=1708},
    @{Row=63; Value=1674},
    @{Row=64; Value=1643},
    @{Row=65; Value=1610},
    @{Row=66; Value=1318},
    @{Row=67; Value=1274},
    @{Row=68; Value=1231},
    @{Row=69; Value=1196},
    @{Row=70; Value=796},
    @{Row=71; Value=746},
    @{Row=72; Value=708},
    @{Row=73; Value=676},
    @{Row=74; Value=309},
    @{Row=75; Value=280},
    @{Row=76; Value=256},
    @{Row=77; Value=237},
    @{Row=78; Value=66},
    @{Row=79; Value=53},
    @{Row=80; Value=44},
    @{Row=81; Value=38},
    @{Row=83; Value=5},
    @{Row=84; Value=5},
    @{Row=85; Value=5},
    @{Row=90; Value=2}
)

foreach ($upd in $bUpdates) {
    $ws.Cells.Item($upd.Row, 2).Value = $upd.Value
}
